$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 34 entirely; this shifts all subsequent rows (35..89) up by one
# row (becoming rows 34..88) and shrinks the used range from A1:P89 to A1:P88,
# matching the commit's automatic map refresh.
$ws.Rows.Item(34).Delete()
